# Update example model: swap the "description" and "type" columns (B and E)
# on the Processes sheet, adjust column widths, move the active
# selection/tab to Processes, and update the data-validation target range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Processes")

# --- Swap columns B (description) and E (type), header + all data rows ---
for ($r = 1; $r -le 9; $r++) {
    $bCell = $ws.Cells.Item($r, 2)
    $eCell = $ws.Cells.Item($r, 5)
    $bVal = $bCell.Value2()
    $eVal = $eCell.Value2()
    $bCell.Value = $eVal
    $eCell.Value = $bVal
}

# --- Column widths: B narrower, E wider, new F column gets B's old width ---
$ws.Columns.Item(2).ColumnWidth = 14.833333333333334
$ws.Columns.Item(5).ColumnWidth = 15.333333333333334
$ws.Columns.Item(6).ColumnWidth = 16.0

# --- Make Processes the active sheet/tab, with E1:E9 selected ---
$ws.Activate()
$ws.Range("E1:E9").Select()
